# "atualizado lista de gastos" — update Fevereiro (sheet4) and Marco (sheet5)
# expense lists.

$wb = $excel.ActiveWorkbook

# --- Fevereiro --------------------------------------------------------
# Item "Capa flip iPhone 5s" (row 5) cost revised from 32.89 to 32.80
$wsFeb = $wb.Worksheets.Item(4)
$wsFeb.Range("B5").Value = 32.8

# --- Marco --------------------------------------------------------------
$wsMar = $wb.Worksheets.Item(5)

# Item "Cartao SD" (row 5) moves from the "Nao gasto" column (C) into the
# "Valor" column (B) and its amount is revised from 109 to 126.
$wsMar.Range("C5").Clear()
$wsMar.Range("B5").Value = 126

# Two new expense rows appended under the existing data.
$wsMar.Range("A7").Value = "Chip Fluke"
$wsMar.Range("B7").Value = 26.97

$wsMar.Range("A8").Value = "Loterias CEF"
$wsMar.Range("B8").Value = 34.5

# --- Restore on-screen selections --------------------------------------
# (Marco is the active tab; select it last so it stays the active sheet.)
$wsFeb.Range("B6").Select()
$wsMar.Range("B9").Select()
